$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 87.98189805347269
$ws.Range("D3").Value = 88.28376452377132
$ws.Range("D4").Value = 87.15383963941923
$ws.Range("D5").Value = 87.85947559951479
$ws.Range("D6").Value = 88.09875854372994
$ws.Range("D7").Value = 88.91898936863468
$ws.Range("D8").Value = 87.6269256467444
$ws.Range("D9").Value = 86.40199837580568
$ws.Range("D10").Value = 85.10970347929863
$ws.Range("D11").Value = 85.82012456093744
$ws.Range("D12").Value = 85.09170501959174
$ws.Range("D13").Value = 85.7680202656022
$ws.Range("D14").Value = 87.35882739828995
$ws.Range("D15").Value = 88.85372779747212
$ws.Range("D16").Value = 88.85603969260613
$ws.Range("D17").Value = 90.6967300113522
$ws.Range("D18").Value = 91.29946932281176
$ws.Range("D19").Value = 92.05869818976858
$ws.Range("D20").Value = 91.20616836197172
$ws.Range("D21").Value = 91.96054185222719
$ws.Range("D22").Value = 92.31096399578379
$ws.Range("D23").Value = 85.41441765623125
$ws.Range("D24").Value = 85.46216809742464
$ws.Range("D25").Value = 84.55102285920174
$ws.Range("D26").Value = 85.24212449847059
$ws.Range("D27").Value = 85.39254559873116
$ws.Range("D28").Value = 86.24122476500217
$ws.Range("D29").Value = 84.21958837190678
$ws.Range("D38").Value = 87.30259251929546
$ws.Range("D39").Value = 88.02800283174703
$ws.Range("D40").Value = 89.13957176843775
$ws.Range("D41").Value = 87.76160329045526
$ws.Range("D42").Value = 88.67195362505535
$ws.Range("D43").Value = 89.1498039836851
$ws.Range("D44").Value = 82.44937319189971
$ws.Range("D45").Value = 84.78873239436619
$ws.Range("D46").Value = 84.47789275634995
$ws.Range("D47").Value = 84.6503178928247
$ws.Range("D48").Value = 85.21897810218978
$ws.Range("D49").Value = 85.06666666666666
$ws.Range("D50").Value = 84.21052631578947
$ws.Range("D59").Value = 87.24954462659382
$ws.Range("D60").Value = 87.87037037037037
$ws.Range("D61").Value = 88.04744525547446
$ws.Range("D62").Value = 88.1740775780511
$ws.Range("D63").Value = 89.76303317535546
$ws.Range("D64").Value = 90.20332717190388
